# Scheduled market-data refresh: updates cached price/profit values
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 115.066666
$ws.Range("I33").Value = 117.875
$ws.Range("K33").Value = 117.875
$ws.Range("M33").Value = 111.125

# Row 46
$ws.Range("H46").Value = 6000
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 18000
$ws.Range("N46").Value = -18238

# Row 60
$ws.Range("H60").Value = 6000
$ws.Range("J60").Value = 6000
$ws.Range("L60").Value = 18000
$ws.Range("N60").Value = -18968

# Row 95
$ws.Range("H95").Value = 38000
$ws.Range("J95").Value = 38000
$ws.Range("L95").Value = 38000
$ws.Range("N95").Value = -43492

# Row 97
$ws.Range("H97").Value = 1197
$ws.Range("J97").Value = 1197
$ws.Range("L97").Value = 3591
$ws.Range("N97").Value = -4583

# Row 129
$ws.Range("H129").Value = 900.6111
$ws.Range("I129").Value = 1325
$ws.Range("K129").Value = 3975
$ws.Range("M129").Value = 1025

# Row 132
$ws.Range("H132").Value = 1238.2142
$ws.Range("I132").Value = 1174.2916
$ws.Range("K132").Value = 3522.8748
$ws.Range("M132").Value = -992.8748000000001

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3238.4465
$ws.Range("I32").Value = 2170.5
$ws.Range("K32").Value = 2170.5
$ws.Range("M32").Value = -1883.5

# Row 61
$ws.Range("H61").Value = 3494.4849
$ws.Range("I61").Value = 2596
$ws.Range("K61").Value = 2596
$ws.Range("M61").Value = -2384

# Row 122
$ws.Range("H122").Value = 5316.5
$ws.Range("I122").Value = 5316.5
$ws.Range("K122").Value = 15949.5
$ws.Range("M122").Value = -13499.5

# Row 136
$ws.Range("H136").Value = 3494.4849
$ws.Range("I136").Value = 2596
$ws.Range("K136").Value = 7788
$ws.Range("M136").Value = -5238

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 73109.21000000001
$ws.Range("I86").Value = 1546.619
$ws.Range("J86").Value = 287797
$ws.Range("K86").Value = 1546.619
$ws.Range("L86").Value = 287797
$ws.Range("M86").Value = -423.6189999999999
$ws.Range("N86").Value = -290043

# Row 89
$ws.Range("H89").Value = 73109.21000000001
$ws.Range("I89").Value = 1546.619
$ws.Range("J89").Value = 287797
$ws.Range("K89").Value = 7733.094999999999
$ws.Range("L89").Value = 1438985
$ws.Range("M89").Value = -2117.094999999999
$ws.Range("N89").Value = -1450217

# Row 94
$ws.Range("H94").Value = 640.84
$ws.Range("I94").Value = 614.9091
$ws.Range("J94").Value = 831
$ws.Range("K94").Value = 614.9091
$ws.Range("L94").Value = 831
$ws.Range("M94").Value = -163.9091
$ws.Range("N94").Value = -1733

# Row 134
$ws.Range("H134").Value = 7075.2095
$ws.Range("I134").Value = 6832
$ws.Range("K134").Value = 20496
$ws.Range("M134").Value = -17961

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1430
$ws.Range("I22").Value = 290
$ws.Range("K22").Value = 290
$ws.Range("M22").Value = 60

# Row 86
$ws.Range("H86").Value = 142859220
$ws.Range("I86").Value = 250001460
$ws.Range("K86").Value = 250001460
$ws.Range("M86").Value = -250000337

# Row 89
$ws.Range("H89").Value = 142859220
$ws.Range("I89").Value = 250001460
$ws.Range("K89").Value = 1250007300
$ws.Range("M89").Value = -1250001684

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 815.55554
$ws.Range("I5").Value = 802
$ws.Range("K5").Value = 2406
$ws.Range("M5").Value = -2294

# Row 56
$ws.Range("H56").Value = 8241.272000000001
$ws.Range("I56").Value = 8241.272000000001
$ws.Range("K56").Value = 8241.272000000001
$ws.Range("M56").Value = -7711.272000000001

# Row 131
$ws.Range("H131").Value = 9743.423000000001
$ws.Range("I131").Value = 618
$ws.Range("J131").Value = 10643.112
$ws.Range("K131").Value = 1854
$ws.Range("L131").Value = 31929.336
$ws.Range("M131").Value = 3186
$ws.Range("N131").Value = -42009.336

# Row 135
$ws.Range("H135").Value = 815.55554
$ws.Range("I135").Value = 802
$ws.Range("K135").Value = 7218
$ws.Range("M135").Value = -4683

# Row 137
$ws.Range("H137").Value = 3378.5
$ws.Range("I137").Value = 1972.4286
$ws.Range("J137").Value = 4135.615
$ws.Range("K137").Value = 5917.2858
$ws.Range("L137").Value = 12406.845
$ws.Range("M137").Value = -817.2857999999997
$ws.Range("N137").Value = -22606.845

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3213.2144
$ws.Range("I102").Value = 3496.7778
$ws.Range("J102").Value = 2702.8
$ws.Range("K102").Value = 3496.7778
$ws.Range("L102").Value = 2702.8
$ws.Range("M102").Value = -1874.7778
$ws.Range("N102").Value = -5946.8

# Row 132
$ws.Range("H132").Value = 1070948.2
$ws.Range("I132").Value = 1540430.2
$ws.Range("J132").Value = 3943.5454
$ws.Range("K132").Value = 4621290.6
$ws.Range("L132").Value = 11830.6362
$ws.Range("M132").Value = -4618760.6
$ws.Range("N132").Value = -16890.6362

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6237.3335
$ws.Range("I7").Value = 3380.8
$ws.Range("J7").Value = 8277.714
$ws.Range("K7").Value = 3380.8
$ws.Range("L7").Value = 8277.714
$ws.Range("M7").Value = -3268.8
$ws.Range("N7").Value = -8501.714

# Row 46
$ws.Range("H46").Value = 2354.2727
$ws.Range("I46").Value = 1150
$ws.Range("K46").Value = 1150
$ws.Range("M46").Value = -962

# Row 82
$ws.Range("H82").Value = 3810
$ws.Range("I82").Value = 1125
$ws.Range("K82").Value = 1125
$ws.Range("M82").Value = -764

# Row 85
$ws.Range("H85").Value = 3810
$ws.Range("I85").Value = 1125
$ws.Range("K85").Value = 1125
$ws.Range("M85").Value = 123

# Row 93
$ws.Range("H93").Value = 460.54544
$ws.Range("I93").Value = 392.55554
$ws.Range("J93").Value = 766.5
$ws.Range("K93").Value = 392.55554
$ws.Range("L93").Value = 766.5
$ws.Range("M93").Value = 855.4444599999999
$ws.Range("N93").Value = -3262.5

# Row 126
$ws.Range("H126").Value = 6237.3335
$ws.Range("I126").Value = 3380.8
$ws.Range("J126").Value = 8277.714
$ws.Range("K126").Value = 10142.4
$ws.Range("L126").Value = 24833.142
$ws.Range("M126").Value = -7672.400000000001
$ws.Range("N126").Value = -29773.142

# Row 132
$ws.Range("H132").Value = 2026.6923
$ws.Range("I132").Value = 1965.091
$ws.Range("K132").Value = 5895.272999999999
$ws.Range("M132").Value = -3365.272999999999

# Row 136
$ws.Range("H136").Value = 3157.9375
$ws.Range("I136").Value = 1503
$ws.Range("K136").Value = 4509
$ws.Range("M136").Value = -1959

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 16837328
$ws.Range("I136").Value = 24156624
$ws.Range("K136").Value = 72469872
$ws.Range("M136").Value = -72467322
